# Generate Report for Handoff
# Updates the localization-status workbook after a handoff event for the
# 547d8ab5-cf9e-4ca6-aad0-cf52a048a451 and 56ffcd2e-04b2-4689-9245-d4efa9dfd71d
# files: status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", timestamps are refreshed, and an "Error Detail"
# message is recorded for the zh-cn / de-de per-file sheets (the handback
# file is stale relative to the latest source revision).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "Overview": rows 4 (547d8ab5...) and 5 (56ffcd2e...)
# Columns E (zh-cn) / F (de-de) hold the consolidated status text,
# column G holds the "Latest HO Xliff Generate Date".
# ---------------------------------------------------------------------
$wsOverview.Range("E4").Value = $readyStatus
$wsOverview.Range("F4").Value = $readyStatus
$wsOverview.Range("G4").Value = "2016-08-17 10:24:41"

$wsOverview.Range("E5").Value = $readyStatus
$wsOverview.Range("F5").Value = $readyStatus
$wsOverview.Range("G5").Value = "2016-08-17 10:24:41"

# ---------------------------------------------------------------------
# Sheet "zh-cn": per-file detail rows 4 (547d8ab5...) and 5 (56ffcd2e...)
# Column C = Status, H = Latest Handoff Datetime, P = Error Detail.
# ---------------------------------------------------------------------
$wsZhCn.Range("C4").Value = $readyStatus
$wsZhCn.Range("H4").Value = "2016-08-17 10:24:36"
$wsZhCn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9dd2012d8b3b7531da99c1cc062eb7359f25154d/e2e/547d8ab5-cf9e-4ca6-aad0-cf52a048a451.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c84e4cee307704ab456fec05cd21158c8244e55/e2e/547d8ab5-cf9e-4ca6-aad0-cf52a048a451.md."

$wsZhCn.Range("C5").Value = $readyStatus
$wsZhCn.Range("H5").Value = "2016-08-17 10:24:36"
$wsZhCn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9dd2012d8b3b7531da99c1cc062eb7359f25154d/e2e/56ffcd2e-04b2-4689-9245-d4efa9dfd71d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c84e4cee307704ab456fec05cd21158c8244e55/e2e/56ffcd2e-04b2-4689-9245-d4efa9dfd71d.md."

# Widen the Error Detail column so the new long messages are readable.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# Sheet "de-de": per-file detail rows 4 (547d8ab5...) and 5 (56ffcd2e...)
# Column C = Status, H = Latest Handoff Datetime, P = Error Detail.
# ---------------------------------------------------------------------
$wsDeDe.Range("C4").Value = $readyStatus
$wsDeDe.Range("H4").Value = "2016-08-17 10:24:41"
$wsDeDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9dd2012d8b3b7531da99c1cc062eb7359f25154d/e2e/547d8ab5-cf9e-4ca6-aad0-cf52a048a451.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c84e4cee307704ab456fec05cd21158c8244e55/e2e/547d8ab5-cf9e-4ca6-aad0-cf52a048a451.md."

$wsDeDe.Range("C5").Value = $readyStatus
$wsDeDe.Range("H5").Value = "2016-08-17 10:24:41"
$wsDeDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9dd2012d8b3b7531da99c1cc062eb7359f25154d/e2e/56ffcd2e-04b2-4689-9245-d4efa9dfd71d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c84e4cee307704ab456fec05cd21158c8244e55/e2e/56ffcd2e-04b2-4689-9245-d4efa9dfd71d.md."

# Widen the Error Detail column so the new long messages are readable.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
